$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-5 with new Name values, keep/update Group as needed
$ws.Range("A2").Value = "SS1381-01-1"
$ws.Range("B2").Value = "young"

$ws.Range("A3").Value = "SS1381-01-2"
$ws.Range("B3").Value = "young"

$ws.Range("A4").Value = "SS1381-01-3A"
$ws.Range("B4").Value = "young"

$ws.Range("A5").Value = "SS1381-01-4A"
$ws.Range("B5").Value = "young"

# Add new rows 6-11
$ws.Range("A6").Value = "SS1381-01-5"
$ws.Range("B6").Value = "young"

$ws.Range("A7").Value = "SS1381-23-4A"
$ws.Range("B7").Value = "old"

$ws.Range("A8").Value = "SS1381-23-13A"
$ws.Range("B8").Value = "old"

$ws.Range("A9").Value = "SS1381-23-14A"
$ws.Range("B9").Value = "old"

$ws.Range("A10").Value = "SS1381-23-15A"
$ws.Range("B10").Value = "old"

$ws.Range("A11").Value = "SS1381-23-16A"
$ws.Range("B11").Value = "old"

# Update selection to match final state
$ws.Range("A5").Select()
